# Suite.xlsx: commit "Checking in Suite.xlsx file"
#
# Substantive change in the "Test Suite" sheet: the Runmode column (C) had
# "N" for rows 3-7 (TSID B Suite .. F Suite); these are flipped to "Y", so
# that every row's Runmode now reads "Y". Since the literal string "N" is
# no longer referenced anywhere in the workbook, it naturally drops out of
# the shared-strings table (uniqueCount 17 -> 16) when the file is saved.
# The active selection also moves from C3 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("C3:C7").Value = "Y"

$ws.Range("C2").Select()
